$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 40 and 41 (F:V) ---
$row40vals = $ws.Range("F40:V40").Value2
$row41vals = $ws.Range("F41:V41").Value2
$ws.Range("F40:V40").Value2 = $row41vals
$ws.Range("F41:V41").Value2 = $row40vals

# --- Append new rows 49-53, copying formatting from row 48 ---
$ws.Range("A48:V48").Copy()
$ws.Range("A49:V49").PasteSpecial(-4122)
$ws.Range("A49").Value2 = 48
$ws.Range("B49").Value2 = 'lebanon'
$ws.Range("C49").Value2 = 'premier-league'
$ws.Range("D49").Value2 = '2023-2024'
$ws.Range("E49").Value2 = 45267.55208333334
$ws.Range("F49").Value2 = 'Al Ahed'
$ws.Range("G49").Value2 = 5
$ws.Range("H49").Value2 = 'Al Ghazieh'
$ws.Range("I49").Value2 = 0
$ws.Range("J49").Value2 = 1.06
$ws.Range("K49").Value2 = '07/12/2023 11:51'
$ws.Range("L49").Value2 = 1.06
$ws.Range("M49").Value2 = '07/12/2023 11:51'
$ws.Range("N49").Value2 = 9.39
$ws.Range("O49").Value2 = '07/12/2023 11:51'
$ws.Range("P49").Value2 = 9.4
$ws.Range("Q49").Value2 = '07/12/2023 11:54'
$ws.Range("R49").Value2 = 24.87
$ws.Range("S49").Value2 = '07/12/2023 11:51'
$ws.Range("T49").Value2 = 24.75
$ws.Range("U49").Value2 = '07/12/2023 11:54'
$ws.Range("V49").Value2 = 'https://www.betexplorer.com/football/lebanon/premier-league/al-ahed-al-ghazieh/GtARTWB3/'

$ws.Range("A48:V48").Copy()
$ws.Range("A50:V50").PasteSpecial(-4122)
$ws.Range("A50").Value2 = 49
$ws.Range("B50").Value2 = 'lebanon'
$ws.Range("C50").Value2 = 'premier-league'
$ws.Range("D50").Value2 = '2023-2024'
$ws.Range("E50").Value2 = 45268.55208333334
$ws.Range("F50").Value2 = 'Tadamon'
$ws.Range("G50").Value2 = 0
$ws.Range("H50").Value2 = 'Bourj FC'
$ws.Range("I50").Value2 = 1
$ws.Range("J50").Value2 = 4.47
$ws.Range("K50").Value2 = '07/12/2023 01:42'
$ws.Range("L50").Value2 = 4.49
$ws.Range("M50").Value2 = '08/12/2023 12:44'
$ws.Range("N50").Value2 = 3.26
$ws.Range("O50").Value2 = '07/12/2023 01:42'
$ws.Range("P50").Value2 = 3.26
$ws.Range("Q50").Value2 = '08/12/2023 11:19'
$ws.Range("R50").Value2 = 1.74
$ws.Range("S50").Value2 = '07/12/2023 01:42'
$ws.Range("T50").Value2 = 1.74
$ws.Range("U50").Value2 = '08/12/2023 12:44'
$ws.Range("V50").Value2 = 'https://www.betexplorer.com/football/lebanon/premier-league/tadamon-bourj/KEGIVhsi/'

$ws.Range("A48:V48").Copy()
$ws.Range("A51:V51").PasteSpecial(-4122)
$ws.Range("A51").Value2 = 50
$ws.Range("B51").Value2 = 'lebanon'
$ws.Range("C51").Value2 = 'premier-league'
$ws.Range("D51").Value2 = '2023-2024'
$ws.Range("E51").Value2 = 45268.625
$ws.Range("F51").Value2 = 'Nejmeh SC'
$ws.Range("G51").Value2 = 1
$ws.Range("H51").Value2 = 'Al Ansar'
$ws.Range("I51").Value2 = 5
$ws.Range("J51").Value2 = 1.88
$ws.Range("K51").Value2 = '07/12/2023 03:12'
$ws.Range("L51").Value2 = 1.94
$ws.Range("M51").Value2 = '08/12/2023 14:15'
$ws.Range("N51").Value2 = 3.16
$ws.Range("O51").Value2 = '07/12/2023 03:12'
$ws.Range("P51").Value2 = 3.11
$ws.Range("Q51").Value2 = '08/12/2023 13:04'
$ws.Range("R51").Value2 = 3.85
$ws.Range("S51").Value2 = '07/12/2023 03:12'
$ws.Range("T51").Value2 = 3.72
$ws.Range("U51").Value2 = '08/12/2023 14:15'
$ws.Range("V51").Value2 = 'https://www.betexplorer.com/football/lebanon/premier-league/nejmeh-sc-al-ansar/fX8VSjR9/'

$ws.Range("A48:V48").Copy()
$ws.Range("A52:V52").PasteSpecial(-4122)
$ws.Range("A52").Value2 = 51
$ws.Range("B52").Value2 = 'lebanon'
$ws.Range("C52").Value2 = 'premier-league'
$ws.Range("D52").Value2 = '2023-2024'
$ws.Range("E52").Value2 = 45270.55208333334
$ws.Range("F52").Value2 = 'Al Ahli Nabatiya'
$ws.Range("G52").Value2 = 0
$ws.Range("H52").Value2 = 'Al Sahel'
$ws.Range("I52").Value2 = 0
$ws.Range("J52").Value2 = 6.29
$ws.Range("K52").Value2 = '10/12/2023 02:42'
$ws.Range("L52").Value2 = 6.34
$ws.Range("M52").Value2 = '10/12/2023 11:19'
$ws.Range("N52").Value2 = 3.66
$ws.Range("O52").Value2 = '10/12/2023 02:42'
$ws.Range("P52").Value2 = 3.4
$ws.Range("Q52").Value2 = '10/12/2023 11:19'
$ws.Range("R52").Value2 = 1.47
$ws.Range("S52").Value2 = '10/12/2023 02:42'
$ws.Range("T52").Value2 = 1.53
$ws.Range("U52").Value2 = '10/12/2023 11:19'
$ws.Range("V52").Value2 = 'https://www.betexplorer.com/football/lebanon/premier-league/al-ahli-nabatiya-al-sahel/h8OW7kkp/'

$ws.Range("A48:V48").Copy()
$ws.Range("A53:V53").PasteSpecial(-4122)
$ws.Range("A53").Value2 = 52
$ws.Range("B53").Value2 = 'lebanon'
$ws.Range("C53").Value2 = 'premier-league'
$ws.Range("D53").Value2 = '2023-2024'
$ws.Range("E53").Value2 = 45270.625
$ws.Range("F53").Value2 = 'Safa'
$ws.Range("G53").Value2 = 2
$ws.Range("H53").Value2 = 'Racing'
$ws.Range("I53").Value2 = 2
$ws.Range("J53").Value2 = 1.96
$ws.Range("K53").Value2 = '10/12/2023 04:12'
$ws.Range("L53").Value2 = 1.93
$ws.Range("M53").Value2 = '10/12/2023 13:05'
$ws.Range("N53").Value2 = 3.26
$ws.Range("O53").Value2 = '10/12/2023 04:12'
$ws.Range("P53").Value2 = 3.31
$ws.Range("Q53").Value2 = '10/12/2023 13:05'
$ws.Range("R53").Value2 = 3.37
$ws.Range("S53").Value2 = '10/12/2023 04:12'
$ws.Range("T53").Value2 = 3.5
$ws.Range("U53").Value2 = '10/12/2023 13:05'
$ws.Range("V53").Value2 = 'https://www.betexplorer.com/football/lebanon/premier-league/safa-racing/b5FMUCdc/'

$excel.CutCopyMode = 0
